# Inserts a new data row at row 38 of Sheet1 (pushing the existing rows
# 38-99 down to 39-100) and populates it with the new "Poroto verde"
# price-point for "Región del Maule" dated 44536.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 38..99 down one position to make room for the new entry.
$ws.Rows.Item(38).Insert()

# Fill in the newly inserted row 38 with the reported data.
$ws.Cells.Item(38, 1).Value  = 10
$ws.Cells.Item(38, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(38, 3).Value  = "La Araucanía"
$ws.Cells.Item(38, 4).Value  = 44536
$ws.Cells.Item(38, 5).Value  = 9
$ws.Cells.Item(38, 6).Value  = 100112031
$ws.Cells.Item(38, 7).Value  = "Poroto verde"
$ws.Cells.Item(38, 8).Value  = "Sin especificar"
$ws.Cells.Item(38, 9).Value  = "Primera"
$ws.Cells.Item(38, 10).Value = 235
$ws.Cells.Item(38, 11).Value = 23000
$ws.Cells.Item(38, 12).Value = 25000
$ws.Cells.Item(38, 13).Value = 24064
$ws.Cells.Item(38, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(38, 15).Value = "Región del Maule"
$ws.Cells.Item(38, 16).Value = 963
$ws.Cells.Item(38, 17).Value = 25
$ws.Cells.Item(38, 18).Value = "Hortaliza"
